$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style of an existing plain-text data cell (no explicit number format), used to
# keep updated cells formatted the same way as before (General / text, no percent fmt).
$refStyle = $ws.Range("B2").Style

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = $refStyle
}

Set-TextValue "D2" "290.97"
Set-TextValue "E2" "-3.31%"
Set-TextValue "D3" "30.70"
Set-TextValue "E3" "-4.80%"
Set-TextValue "D4" "4.953"
Set-TextValue "E4" "-0.03%"
Set-TextValue "D5" "0.07216"
Set-TextValue "E5" "-5.56%"
Set-TextValue "D6" "1.844"
Set-TextValue "E6" "-3.71%"
Set-TextValue "D7" "7.698"
Set-TextValue "E7" "-1.72%"
Set-TextValue "D8" "3.772"
Set-TextValue "E8" "-0.76%"
Set-TextValue "D9" "0.8971"
Set-TextValue "E9" "-2.21%"
Set-TextValue "D10" "0.1660"
Set-TextValue "E10" "-5.07%"
Set-TextValue "D11" "0.07726"
Set-TextValue "E11" "-0.31%"
Set-TextValue "D12" "0.07995"
Set-TextValue "E12" "-6.69%"
Set-TextValue "D13" "0.03042"
Set-TextValue "E13" "-6.07%"
Set-TextValue "D14" "0.1001"
Set-TextValue "E14" "0.06%"
Set-TextValue "D15" "0.001498"
Set-TextValue "E15" "-0.89%"
Set-TextValue "D16" "0.005712"
Set-TextValue "E16" "-3.83%"
Set-TextValue "D18" "3.465"
Set-TextValue "E18" "0.01%"
Set-TextValue "E19" "-3.26%"
Set-TextValue "E20" "-0.94%"
Set-TextValue "D21" "0.1297"
Set-TextValue "E21" "-2.20%"
Set-TextValue "D22" "4.052"
Set-TextValue "E22" "-5.55%"
Set-TextValue "D23" "0.2322"
Set-TextValue "E23" "16.60%"
Set-TextValue "D24" "0.04508"
Set-TextValue "E24" "-0.40%"
Set-TextValue "D25" "0.001217"
Set-TextValue "D26" "0.004647"
Set-TextValue "E26" "5.91%"
Set-TextValue "E27" "-0.09%"
Set-TextValue "D39" "0.01578"
Set-TextValue "E39" "-7.12%"
Set-TextValue "E40" "-6.29%"
Set-TextValue "D41" "0.007258"
Set-TextValue "E41" "-3.18%"
Set-TextValue "D43" "0.1303"
Set-TextValue "E43" "-3.43%"
Set-TextValue "E44" "-11.67%"
Set-TextValue "D45" "0.009194"
Set-TextValue "E45" "-12.77%"
Set-TextValue "D46" "0.00005969"
Set-TextValue "E46" "-4.39%"
Set-TextValue "E47" "-0.08%"
Set-TextValue "D48" "2.247"
Set-TextValue "E48" "173.93%"
Set-TextValue "E50" "-0.08%"
Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "-0.08%"
